# Update the "刷题记录表" (LeetCode practice log) worksheet:
#  - Add row 12: problem 154 "Find Minimum in Rotated Sorted Array II"
#  - Add row 13: problem 278 "First Bad Version"
#  - Update the selection / scroll position to the new last cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 -----------------------------------------------------------
# Duplicate the last existing data row (row 11) so the new row picks up
# the same cell styles (number formats, alignment, wrap) and row height
# instead of creating brand-new style/number-format entries.
$ws.Rows.Item(11).Copy()
$ws.Rows.Item(12).Insert(-4121)
$ws.Rows.Item(12).RowHeight = 34

$ws.Range("A12").Value = 154
$ws.Range("B12").Value = "Find Minimum in Rotated Sorted Array II"
$ws.Range("C12").Value = "#array  #binary-search #重点 "
$ws.Range("D12").Value = "hard"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 21
$ws.Range("H12").Value = 45832
$ws.Range("I12").Value = 45832

# --- Row 13 -----------------------------------------------------------
$ws.Rows.Item(12).Copy()
$ws.Rows.Item(13).Insert(-4121)
$ws.Rows.Item(13).RowHeight = 34

$ws.Range("A13").Value = 278
$ws.Range("B13").Value = "First Bad Version"
$ws.Range("C13").Value = "#binary-search #重点"
$ws.Range("D13").Value = "easy"
$ws.Range("E13").Value = 6
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = 45832
$ws.Range("I13").Value = 45832

# --- View state ---------------------------------------------------------
[void]$ws.Range("I13").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
